$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers share the bold/bordered/centered header style used by the
# rest of row 1 (copy format from the last existing header cell, AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (constant across the roster) for every player row.
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
